$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.62197818638615
$ws.Cells.Item(2, 3).Value = 6.668666515371771
$ws.Cells.Item(2, 5).Value = 14.24885395091472
$ws.Cells.Item(2, 6).Value = 41.08018402194973
$ws.Cells.Item(2, 7).Value = 42.14033072464268
$ws.Cells.Item(2, 8).Value = 17.69607241463209
$ws.Cells.Item(2, 9).Value = 27.92163377542688
$ws.Cells.Item(2, 10).Value = 8.612173447054335
$ws.Cells.Item(2, 11).Value = 11.51466277271392
$ws.Cells.Item(2, 12).Value = 12.38526435680083
$ws.Cells.Item(2, 14).Value = 20.79034900346628
$ws.Cells.Item(3, 2).Value = 14.42250144070761
$ws.Cells.Item(3, 3).Value = 6.60589082936156
$ws.Cells.Item(3, 5).Value = 14.24609221230732
$ws.Cells.Item(3, 6).Value = 41.08447669867266
$ws.Cells.Item(3, 7).Value = 42.20980034299974
$ws.Cells.Item(3, 8).Value = 17.7448965137716
$ws.Cells.Item(3, 9).Value = 28.00444543028237
$ws.Cells.Item(3, 10).Value = 8.615619158639131
$ws.Cells.Item(3, 11).Value = 11.38142939262293
$ws.Cells.Item(3, 12).Value = 12.37054441348269
$ws.Cells.Item(3, 14).Value = 20.85428197250905
$ws.Cells.Item(4, 2).Value = 14.30165741991432
$ws.Cells.Item(4, 3).Value = 6.566280698436031
$ws.Cells.Item(4, 5).Value = 14.24631360354817
$ws.Cells.Item(4, 6).Value = 41.09602000046888
$ws.Cells.Item(4, 7).Value = 42.26325840878285
$ws.Cells.Item(4, 8).Value = 17.77756375042323
$ws.Cells.Item(4, 9).Value = 28.06005358345505
$ws.Cells.Item(4, 10).Value = 8.617893374407734
$ws.Cells.Item(4, 11).Value = 11.30115439261113
$ws.Cells.Item(4, 12).Value = 12.36333272182989
$ws.Cells.Item(4, 14).Value = 20.89535373743827
$ws.Cells.Item(5, 2).Value = 14.25288352610617
$ws.Cells.Item(5, 3).Value = 6.549876193732726
$ws.Cells.Item(5, 5).Value = 14.24688756688911
$ws.Cells.Item(5, 6).Value = 41.10296508661078
$ws.Cells.Item(5, 7).Value = 42.28775221370141
$ws.Cells.Item(5, 8).Value = 17.79155184797261
$ws.Cells.Item(5, 9).Value = 28.08391050590268
$ws.Cells.Item(5, 10).Value = 8.618860137065939
$ws.Cells.Item(5, 11).Value = 11.26886304685481
$ws.Cells.Item(5, 12).Value = 12.36085580775575
$ws.Cells.Item(5, 14).Value = 20.91254897677234
$ws.Cells.Item(6, 2).Value = 14.24481484072886
$ws.Cells.Item(6, 3).Value = 6.547136495154901
$ws.Cells.Item(6, 5).Value = 14.24701213569437
$ws.Cells.Item(6, 6).Value = 41.10425368255457
$ws.Cells.Item(6, 7).Value = 42.29198278523762
$ws.Cells.Item(6, 8).Value = 17.79391537817471
$ws.Cells.Item(6, 9).Value = 28.08794413959259
$ws.Cells.Item(6, 10).Value = 8.619023087246601
$ws.Cells.Item(6, 11).Value = 11.26352759266139
$ws.Cells.Item(6, 12).Value = 12.36047248847764
$ws.Cells.Item(6, 14).Value = 20.91543194461996
$ws.Cells.Item(7, 2).Value = 14.30099765248348
$ws.Cells.Item(7, 3).Value = 6.566060520105771
$ws.Cells.Item(7, 5).Value = 14.24631938346064
$ws.Cells.Item(7, 6).Value = 41.09610458948852
$ws.Cells.Item(7, 7).Value = 42.26357778187889
$ws.Cells.Item(7, 8).Value = 17.77774966232123
$ws.Cells.Item(7, 9).Value = 28.06037048467331
$ws.Cells.Item(7, 10).Value = 8.61790625036514
$ws.Cells.Item(7, 11).Value = 11.30071714656846
$ws.Cells.Item(7, 12).Value = 12.3632974437436
$ws.Cells.Item(7, 14).Value = 20.89558378159402
$ws.Cells.Item(8, 2).Value = 14.55289517472322
$ws.Cells.Item(8, 3).Value = 6.647245532635919
$ws.Cells.Item(8, 5).Value = 14.24750485541748
$ws.Cells.Item(8, 6).Value = 41.07981642682634
$ws.Cells.Item(8, 7).Value = 42.16203743164962
$ws.Cells.Item(8, 8).Value = 17.71234870834529
$ws.Cells.Item(8, 9).Value = 27.94919830653723
$ws.Cells.Item(8, 10).Value = 8.613328706453684
$ws.Cells.Item(8, 11).Value = 11.46842824789417
$ws.Cells.Item(8, 12).Value = 12.37981150532107
$ws.Cells.Item(8, 14).Value = 20.81201685750911
$ws.Cells.Item(9, 2).Value = 15.05719850196868
$ws.Cells.Item(9, 3).Value = 6.797785888345071
$ws.Cells.Item(9, 5).Value = 14.26496229142338
$ws.Cells.Item(9, 6).Value = 41.11845519372195
$ws.Cells.Item(9, 7).Value = 42.04892794052617
$ws.Cells.Item(9, 8).Value = 17.60544492264519
$ws.Cells.Item(9, 9).Value = 27.7690211000544
$ws.Cells.Item(9, 10).Value = 8.605603990758533
$ws.Cells.Item(9, 11).Value = 11.8078255885156
$ws.Cells.Item(9, 12).Value = 12.42656185991925
$ws.Cells.Item(9, 14).Value = 20.66249414761538
$ws.Cells.Item(10, 2).Value = 15.43047608458542
$ws.Cells.Item(10, 3).Value = 6.902838598842028
$ws.Cells.Item(10, 5).Value = 14.28690028563671
$ws.Cells.Item(10, 6).Value = 41.18967300508471
$ws.Cells.Item(10, 7).Value = 42.01859686508089
$ws.Cells.Item(10, 8).Value = 17.53993127961636
$ws.Cells.Item(10, 9).Value = 27.65977625738174
$ws.Cells.Item(10, 10).Value = 8.60068351583671
$ws.Cells.Item(10, 11).Value = 12.06138552527485
$ws.Cells.Item(10, 12).Value = 12.46948568780927
$ws.Cells.Item(10, 14).Value = 20.56129969232529
$ws.Cells.Item(11, 2).Value = 15.60015404594304
$ws.Cells.Item(11, 3).Value = 6.949370294570747
$ws.Cells.Item(11, 5).Value = 14.29883147062403
$ws.Cells.Item(11, 6).Value = 41.23130482848281
$ws.Cells.Item(11, 7).Value = 42.01630597515596
$ws.Cells.Item(11, 8).Value = 17.51295942821999
$ws.Cells.Item(11, 9).Value = 27.61511472023393
$ws.Cells.Item(11, 10).Value = 8.598607217893555
$ws.Cells.Item(11, 11).Value = 12.17718098229851
$ws.Cells.Item(11, 12).Value = 12.4908312590713
$ws.Cells.Item(11, 14).Value = 20.51712503956091
$ws.Cells.Item(12, 2).Value = 15.66432847599805
$ws.Cells.Item(12, 3).Value = 6.966805315208561
$ws.Cells.Item(12, 5).Value = 14.30362746077517
$ws.Cells.Item(12, 6).Value = 41.24838976016871
$ws.Cells.Item(12, 7).Value = 42.01709507871148
$ws.Cells.Item(12, 8).Value = 17.50315318197936
$ws.Cells.Item(12, 9).Value = 27.59892763583318
$ws.Cells.Item(12, 10).Value = 8.597844139919486
$ws.Cells.Item(12, 11).Value = 12.22105567503185
$ws.Cells.Item(12, 12).Value = 12.4991717338001
$ws.Cells.Item(12, 14).Value = 20.50066319465375
$ws.Cells.Item(13, 2).Value = 15.65051182000977
$ws.Cells.Item(13, 3).Value = 6.963058696488265
$ws.Cells.Item(13, 5).Value = 14.30258224343864
$ws.Cells.Item(13, 6).Value = 41.2446516546607
$ws.Cells.Item(13, 7).Value = 42.01685143519187
$ws.Cells.Item(13, 8).Value = 17.50524700731764
$ws.Cells.Item(13, 9).Value = 27.60238153708353
$ws.Cells.Item(13, 10).Value = 8.598007453972588
$ws.Cells.Item(13, 11).Value = 12.21160595912598
$ws.Cells.Item(13, 12).Value = 12.49736408689625
$ws.Cells.Item(13, 14).Value = 20.50419673019097
$ws.Cells.Item(14, 2).Value = 15.60543563738419
$ws.Cells.Item(14, 3).Value = 6.950808428930991
$ws.Cells.Item(14, 5).Value = 14.29922048605809
$ws.Cells.Item(14, 6).Value = 41.23268401462492
$ws.Cells.Item(14, 7).Value = 42.01633768749965
$ws.Cells.Item(14, 8).Value = 17.51214449533324
$ws.Cells.Item(14, 9).Value = 27.61376845620664
$ws.Cells.Item(14, 10).Value = 8.598543975382146
$ws.Cells.Item(14, 11).Value = 12.18079029139913
$ws.Cells.Item(14, 12).Value = 12.49151230204018
$ws.Cells.Item(14, 14).Value = 20.51576538680466
$ws.Cells.Item(15, 2).Value = 15.57781312764498
$ws.Cells.Item(15, 3).Value = 6.943280487948329
$ws.Cells.Item(15, 5).Value = 14.29719742112084
$ws.Cells.Item(15, 6).Value = 41.22552511446032
$ws.Cells.Item(15, 7).Value = 42.01623877850556
$ws.Cells.Item(15, 8).Value = 17.51642247308213
$ws.Cells.Item(15, 9).Value = 27.62083776274978
$ws.Cells.Item(15, 10).Value = 8.598875623878415
$ws.Cells.Item(15, 11).Value = 12.16191696967533
$ws.Cells.Item(15, 12).Value = 12.48796130192181
$ws.Cells.Item(15, 14).Value = 20.52288614608691
$ws.Cells.Item(16, 2).Value = 15.41937944337903
$ws.Cells.Item(16, 3).Value = 6.899771935998862
$ws.Cells.Item(16, 5).Value = 14.28615959569456
$ws.Cells.Item(16, 6).Value = 41.18713743022381
$ws.Cells.Item(16, 7).Value = 42.01897832201173
$ws.Cells.Item(16, 8).Value = 17.54175094490647
$ws.Cells.Item(16, 9).Value = 27.66279641445687
$ws.Cells.Item(16, 10).Value = 8.600822456402712
$ws.Cells.Item(16, 11).Value = 12.05382370344603
$ws.Cells.Item(16, 12).Value = 12.46812695644843
$ws.Cells.Item(16, 14).Value = 20.56422391775391
$ws.Cells.Item(17, 2).Value = 15.32210915239132
$ws.Cells.Item(17, 3).Value = 6.872755418603979
$ws.Cells.Item(17, 5).Value = 14.27988619553131
$ws.Cells.Item(17, 6).Value = 41.16594768108561
$ws.Cells.Item(17, 7).Value = 42.02360790310534
$ws.Cells.Item(17, 8).Value = 17.55801440206164
$ws.Cells.Item(17, 9).Value = 27.68982709556727
$ws.Cells.Item(17, 10).Value = 8.602058186417196
$ws.Cells.Item(17, 11).Value = 11.98759908206181
$ws.Cells.Item(17, 12).Value = 12.45642225465475
$ws.Cells.Item(17, 14).Value = 20.59005862256593
$ws.Cells.Item(18, 2).Value = 15.26615334136598
$ws.Cells.Item(18, 3).Value = 6.857098162694896
$ws.Cells.Item(18, 5).Value = 14.27646175783634
$ws.Cells.Item(18, 6).Value = 41.15462965298488
$ws.Cells.Item(18, 7).Value = 42.02735375867847
$ws.Cells.Item(18, 8).Value = 17.56763512422565
$ws.Cells.Item(18, 9).Value = 27.70584822281917
$ws.Cells.Item(18, 10).Value = 8.602784206474999
$ws.Cells.Item(18, 11).Value = 11.94955307691769
$ws.Cells.Item(18, 12).Value = 12.44986153126043
$ws.Cells.Item(18, 14).Value = 20.60509313123553
$ws.Cells.Item(19, 2).Value = 15.24720807542997
$ws.Cells.Item(19, 3).Value = 6.851776752392054
$ws.Cells.Item(19, 5).Value = 14.27533395925766
$ws.Cells.Item(19, 6).Value = 41.1509471688159
$ws.Cells.Item(19, 7).Value = 42.02880797526834
$ws.Cells.Item(19, 8).Value = 17.57093828045798
$ws.Cells.Item(19, 9).Value = 27.71135403491469
$ws.Cells.Item(19, 10).Value = 8.603032649475043
$ws.Cells.Item(19, 11).Value = 11.93668022658101
$ws.Cells.Item(19, 12).Value = 12.44766976559531
$ws.Cells.Item(19, 14).Value = 20.61021366316991
$ws.Cells.Item(20, 2).Value = 15.33246506588702
$ws.Cells.Item(20, 3).Value = 6.875643627419007
$ws.Cells.Item(20, 5).Value = 14.28053500035345
$ws.Cells.Item(20, 6).Value = 41.16811340389025
$ws.Cells.Item(20, 7).Value = 42.02300296735847
$ws.Cells.Item(20, 8).Value = 17.55625555393254
$ws.Cells.Item(20, 9).Value = 27.68690058774298
$ws.Cells.Item(20, 10).Value = 8.601925062360044
$ws.Cells.Item(20, 11).Value = 11.99464447156961
$ws.Cells.Item(20, 12).Value = 12.45765051931281
$ws.Cells.Item(20, 14).Value = 20.58729036389348
$ws.Cells.Item(21, 2).Value = 15.61867822086843
$ws.Cells.Item(21, 3).Value = 6.954411701461664
$ws.Cells.Item(21, 5).Value = 14.30020039579207
$ws.Cells.Item(21, 6).Value = 41.23616345024799
$ws.Cells.Item(21, 7).Value = 42.01644361835403
$ws.Cells.Item(21, 8).Value = 17.51010747612489
$ws.Cells.Item(21, 9).Value = 27.61040414817535
$ws.Cells.Item(21, 10).Value = 8.598385758232912
$ws.Cells.Item(21, 11).Value = 12.18984121786329
$ws.Cells.Item(21, 12).Value = 12.49322416317989
$ws.Cells.Item(21, 14).Value = 20.51236017866488
$ws.Cells.Item(22, 2).Value = 15.80524820052495
$ws.Cells.Item(22, 3).Value = 7.004807729705345
$ws.Cells.Item(22, 5).Value = 14.31467135534104
$ws.Cells.Item(22, 6).Value = 41.28832664451932
$ws.Cells.Item(22, 7).Value = 42.02181273168834
$ws.Cells.Item(22, 8).Value = 17.48232188473697
$ws.Cells.Item(22, 9).Value = 27.56463730643704
$ws.Cells.Item(22, 10).Value = 8.596207615703239
$ws.Cells.Item(22, 11).Value = 12.31754435865801
$ws.Cells.Item(22, 12).Value = 12.51797160719933
$ws.Cells.Item(22, 14).Value = 20.46493966741555
$ws.Cells.Item(23, 2).Value = 15.70573637509996
$ws.Cells.Item(23, 3).Value = 6.978011082655907
$ws.Cells.Item(23, 5).Value = 14.30680078131244
$ws.Cells.Item(23, 6).Value = 41.25978564596676
$ws.Cells.Item(23, 7).Value = 42.01806329055312
$ws.Cells.Item(23, 8).Value = 17.49693416306644
$ws.Cells.Item(23, 9).Value = 27.5886766475293
$ws.Cells.Item(23, 10).Value = 8.597357822570569
$ws.Cells.Item(23, 11).Value = 12.24938758636231
$ws.Cells.Item(23, 12).Value = 12.50462782704808
$ws.Cells.Item(23, 14).Value = 20.49010739651247
$ws.Cells.Item(24, 2).Value = 15.32778325810436
$ws.Cells.Item(24, 3).Value = 6.874338257547614
$ws.Cells.Item(24, 5).Value = 14.28024110789499
$ws.Cells.Item(24, 6).Value = 41.1671315879473
$ws.Cells.Item(24, 7).Value = 42.02327308155861
$ws.Cells.Item(24, 8).Value = 17.55704988613548
$ws.Cells.Item(24, 9).Value = 27.68822216424612
$ws.Cells.Item(24, 10).Value = 8.60198519919305
$ws.Cells.Item(24, 11).Value = 11.99145916257701
$ws.Cells.Item(24, 12).Value = 12.45709469579467
$ws.Cells.Item(24, 14).Value = 20.58854132724777
$ws.Cells.Item(25, 2).Value = 14.92003170859306
$ws.Cells.Item(25, 3).Value = 6.758016259881093
$ws.Cells.Item(25, 5).Value = 14.25863053140201
$ws.Cells.Item(25, 6).Value = 41.10046183095333
$ws.Cells.Item(25, 7).Value = 42.07027989193582
$ws.Cells.Item(25, 8).Value = 17.6320782840056
$ws.Cells.Item(25, 9).Value = 27.81370649447678
$ws.Cells.Item(25, 10).Value = 8.607560569997975
$ws.Cells.Item(25, 11).Value = 11.71511120317633
$ws.Cells.Item(25, 12).Value = 12.41239437535366
$ws.Cells.Item(25, 14).Value = 20.70141682264604
